$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column "entradas_vendidas" (column I) ---

# Header cell I1: same text + same formatting (bold/border/centered) as the other headers
$ws.Range("I1").Value = "entradas_vendidas"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells I2:I100 -> entradas_vendidas per play
$entradasVendidas = @(
    ,@(2, 1517)
    ,@(3, 3206)
    ,@(4, 2801)
    ,@(5, 2202)
    ,@(6, 27)
    ,@(7, 617)
    ,@(8, 66)
    ,@(9, 3046)
    ,@(10, 24)
    ,@(11, 808)
    ,@(12, 3019)
    ,@(13, 20)
    ,@(14, 1372)
    ,@(15, 22)
    ,@(16, 2811)
    ,@(17, 20)
    ,@(18, 1631)
    ,@(19, 20)
    ,@(20, 2508)
    ,@(21, 1897)
    ,@(22, 1585)
    ,@(23, 29)
    ,@(24, 36)
    ,@(25, 41)
    ,@(26, 20)
    ,@(27, 1481)
    ,@(28, 2308)
    ,@(29, 2368)
    ,@(30, 2012)
    ,@(31, 28)
    ,@(32, 63)
    ,@(33, 3301)
    ,@(34, 31)
    ,@(35, 2392)
    ,@(36, 367)
    ,@(37, 949)
    ,@(38, 1800)
    ,@(39, 20)
    ,@(40, 100)
    ,@(41, 3001)
    ,@(42, 26)
    ,@(43, 2411)
    ,@(44, 2631)
    ,@(45, 2164)
    ,@(46, 23)
    ,@(47, 21)
    ,@(48, 1072)
    ,@(49, 61)
    ,@(50, 196)
    ,@(51, 966)
    ,@(52, 20)
    ,@(53, 1947)
    ,@(54, 1523)
    ,@(55, 94)
    ,@(56, 872)
    ,@(57, 20)
    ,@(58, 20)
    ,@(59, 20)
    ,@(60, 2415)
    ,@(61, 1430)
    ,@(62, 1812)
    ,@(63, 20)
    ,@(64, 912)
    ,@(65, 3024)
    ,@(66, 97)
    ,@(67, 2469)
    ,@(68, 20)
    ,@(69, 3106)
    ,@(70, 20)
    ,@(71, 1631)
    ,@(72, 20)
    ,@(73, 1632)
    ,@(74, 1196)
    ,@(75, 1303)
    ,@(76, 1811)
    ,@(77, 20)
    ,@(78, 20)
    ,@(79, 1077)
    ,@(80, 761)
    ,@(81, 59)
    ,@(82, 1499)
    ,@(83, 2255)
    ,@(84, 32)
    ,@(85, 2208)
    ,@(86, 1034)
    ,@(87, 977)
    ,@(88, 22)
    ,@(89, 24)
    ,@(90, 1672)
    ,@(91, 48)
    ,@(92, 1101)
    ,@(93, 2414)
    ,@(94, 2523)
    ,@(95, 26)
    ,@(96, 84)
    ,@(97, 2611)
    ,@(98, 1043)
    ,@(99, 20)
    ,@(100, 175)
)

foreach ($pair in $entradasVendidas) {
    $row = $pair[0]
    $value = $pair[1]
    $ws.Cells.Item($row, 9).Value = $value
}

# Widen column I to fit the new header/data
$ws.Columns.Item(9).ColumnWidth = 19.5

# Move the active selection to the new column header (I1)
$null = $ws.Range("I1").Select()

Write-Host "done"
